# ArmyOfSoldiers, Hustle, Replenishment and StrategicStrikes images
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to match its purpose
$ws.Name = "ListOfCards"

# Mark the IMG column ("C") as done ("v") for the cards whose artwork
# was finished.
$ws.Range("C29").Value = "v"   # Hustle
$ws.Range("C41").Value = "v"   # Replenishment
$ws.Range("C46").Value = "v"   # StrategicStrikes
$ws.Range("C54").Value = "v"   # ArmyOfSoldiers
$ws.Range("C55").Value = "v"   # CommandingPresence

# Leave the view scrolled down to where the work happened, with the
# last-touched cell selected.
try {
    $excel.ActiveWindow.ScrollRow = 40
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$ws.Range("C54").Select()
